$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "scenario2" rows exercise the Add/Edit/Delete source-data flows ---

# Row 9 (source_startdate): the previously-entered real dates are re-entered as
# quote-prefixed text (so they round-trip as shared strings, matching how the
# automation captured the "Add value" / "Edit value" results for this run).
$ws.Range("C9").Value = "'11/01/2022"
$ws.Range("D9").Value = "'11/03/2022"

# Row 10 (source_enddate): same treatment as row 9.
$ws.Range("C10").Value = "'11/05/2022"
$ws.Range("D10").Value = "'11/10/2022"

# Row 7 (source_name_text_box) - "Edit_source_value" column now reflects the
# updated/edited automation name.
$ws.Range("D7").Value = "Automation_Test_Update"

# Row 8 (source_abbreviation_text_box) - "Edit_source_value" column now
# reflects the updated/edited abbreviation.
$ws.Range("D8").Value = "AUT_UPDT"

# Leave the same cell selected/active as when the sheet was last saved.
[void]$ws.Range("D9").Select()
